$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 550.4
$ws.Range("I41").Value = 666.6667
$ws.Range("J41").Value = 376
$ws.Range("K41").Value = 666.6667
$ws.Range("L41").Value = 376
$ws.Range("M41").Value = -226.6667
$ws.Range("N41").Value = -1256

$ws.Range("H49").Value = 860

$ws.Range("H98").Value = 3192.8333
$ws.Range("I98").Value = 678.63635
$ws.Range("K98").Value = 678.63635
$ws.Range("M98").Value = 819.36365

$ws.Range("H112").Value = 1529.5834
$ws.Range("I112").Value = 1113.75
$ws.Range("K112").Value = 3341.25
$ws.Range("M112").Value = -2233.25

$ws.Range("H122").Value = 3192.8333
$ws.Range("I122").Value = 678.63635
$ws.Range("K122").Value = 2035.90905
$ws.Range("M122").Value = 414.09095

$ws.Range("H132").Value = 1019.2857
$ws.Range("I132").Value = 1055.8889
$ws.Range("K132").Value = 3167.6667
$ws.Range("M132").Value = -637.6666999999998

$ws.Range("H138").Value = 2907.5957
$ws.Range("I138").Value = 3362.1765
$ws.Range("J138").Value = 2807.234
$ws.Range("K138").Value = 10086.5295
$ws.Range("L138").Value = 8421.701999999999
$ws.Range("M138").Value = -4946.529500000001
$ws.Range("N138").Value = -18701.702

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1525.5385
$ws.Range("I45").Value = 1403.091
$ws.Range("K45").Value = 1403.091
$ws.Range("M45").Value = -1026.091

$ws.Range("H53").Value = 19039
$ws.Range("I53").Value = 19039
$ws.Range("K53").Value = 19039
$ws.Range("M53").Value = -18357

$ws.Range("H61").Value = 4492.25
$ws.Range("I61").Value = 1987.4286
$ws.Range("J61").Value = 7999
$ws.Range("K61").Value = 1987.4286
$ws.Range("L61").Value = 7999
$ws.Range("M61").Value = -1775.4286
$ws.Range("N61").Value = -8423

$ws.Range("H132").Value = 1881.8462
$ws.Range("I132").Value = 1622.8182
$ws.Range("J132").Value = 3306.5
$ws.Range("K132").Value = 4868.4546
$ws.Range("L132").Value = 9919.5
$ws.Range("M132").Value = -2338.4546
$ws.Range("N132").Value = -14979.5

$ws.Range("H136").Value = 4492.25
$ws.Range("I136").Value = 1987.4286
$ws.Range("J136").Value = 7999
$ws.Range("K136").Value = 5962.2858
$ws.Range("L136").Value = 23997
$ws.Range("M136").Value = -3412.2858
$ws.Range("N136").Value = -29097

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2308.3572
$ws.Range("I86").Value = 2361.4167
$ws.Range("J86").Value = 1990
$ws.Range("K86").Value = 2361.4167
$ws.Range("L86").Value = 1990
$ws.Range("M86").Value = -1238.4167
$ws.Range("N86").Value = -4236

$ws.Range("H89").Value = 2308.3572
$ws.Range("I89").Value = 2361.4167
$ws.Range("J89").Value = 1990
$ws.Range("K89").Value = 11807.0835
$ws.Range("L89").Value = 9950
$ws.Range("M89").Value = -6191.083500000001
$ws.Range("N89").Value = -21182

$ws.Range("H94").Value = 10358.286
$ws.Range("I94").Value = 10358.286
$ws.Range("K94").Value = 10358.286
$ws.Range("M94").Value = -9907.286

$ws.Range("H105").Value = 3472.5
$ws.Range("I105").Value = 2968.5
$ws.Range("K105").Value = 2968.5
$ws.Range("M105").Value = -1221.5

$ws.Range("H107").Value = 2543.4
$ws.Range("I107").Value = 2350
$ws.Range("J107").Value = 4284
$ws.Range("K107").Value = 2350
$ws.Range("L107").Value = 4284
$ws.Range("M107").Value = -430
$ws.Range("N107").Value = -8124

$ws.Range("H134").Value = 485.75
$ws.Range("I134").Value = 485.75
$ws.Range("K134").Value = 1457.25
$ws.Range("M134").Value = 1077.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2707.4119
$ws.Range("I31").Value = 2345
$ws.Range("K31").Value = 2345
$ws.Range("M31").Value = -2050

$ws.Range("H34").Value = 2707.4119
$ws.Range("I34").Value = 2345
$ws.Range("K34").Value = 2345
$ws.Range("M34").Value = -2143

$ws.Range("H130").Value = 39258.5
$ws.Range("J130").Value = 39258.5
$ws.Range("L130").Value = 39258.5
$ws.Range("N130").Value = -49298.5

$ws.Range("H134").Value = 1708.0333
$ws.Range("I134").Value = 1708.0333
$ws.Range("K134").Value = 5124.0999
$ws.Range("M134").Value = -2589.0999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H19").Value = 1000
$ws.Range("I19").Value = 1000
$ws.Range("J19").Value = 1000
$ws.Range("K19").Value = 3000
$ws.Range("L19").Value = 3000
$ws.Range("M19").Value = -2826
$ws.Range("N19").Value = -3348

$ws.Range("H40").Value = 55.375
$ws.Range("I40").Value = 55.375
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 221.5
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -152.5
$ws.Range("N40").ClearContents()

$ws.Range("H86").Value = 3752
$ws.Range("I86").Value = 652.8
$ws.Range("K86").Value = 1958.4
$ws.Range("M86").Value = -772.3999999999999

$ws.Range("H89").Value = 3752
$ws.Range("I89").Value = 652.8
$ws.Range("K89").Value = 5875.2
$ws.Range("M89").Value = 52.80000000000018

$ws.Range("H122").Value = 496.66666
$ws.Range("J122").Value = 513
$ws.Range("L122").Value = 4617
$ws.Range("N122").Value = -9517

$ws.Range("H127").Value = 1991.4
$ws.Range("J127").Value = 1991.4
$ws.Range("L127").Value = 5974.200000000001
$ws.Range("N127").Value = -15894.2

$ws.Range("H131").Value = 2954.0908
$ws.Range("I131").Value = 1998.3334
$ws.Range("J131").Value = 3312.5
$ws.Range("K131").Value = 5995.0002
$ws.Range("L131").Value = 9937.5
$ws.Range("M131").Value = -955.0002000000004
$ws.Range("N131").Value = -20017.5

$ws.Range("H133").Value = 14666.5
$ws.Range("I133").Value = 3999.5
$ws.Range("K133").Value = 11998.5
$ws.Range("M133").Value = -6938.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3124.2666
$ws.Range("I80").Value = 2201.25
$ws.Range("J80").Value = 3459.9092
$ws.Range("K80").Value = 2201.25
$ws.Range("L80").Value = 3459.9092
$ws.Range("M80").Value = -1203.25
$ws.Range("N80").Value = -5455.9092

$ws.Range("H83").Value = 3124.2666
$ws.Range("I83").Value = 2201.25
$ws.Range("J83").Value = 3459.9092
$ws.Range("K83").Value = 11006.25
$ws.Range("L83").Value = 17299.546
$ws.Range("M83").Value = -6014.25
$ws.Range("N83").Value = -27283.546

$ws.Range("H136").Value = 22475.555
$ws.Range("J136").Value = 22475.555
$ws.Range("L136").Value = 67426.66500000001
$ws.Range("N136").Value = -72526.66500000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6069.0645
$ws.Range("I7").Value = 1941.2222
$ws.Range("K7").Value = 1941.2222
$ws.Range("M7").Value = -1829.2222

$ws.Range("H22").Value = 1068.6
$ws.Range("I22").Value = 864.44446
$ws.Range("J22").Value = 1374.8334
$ws.Range("K22").Value = 864.44446
$ws.Range("L22").Value = 1374.8334
$ws.Range("M22").Value = -569.44446
$ws.Range("N22").Value = -1964.8334

$ws.Range("H27").Value = 1068.6
$ws.Range("I27").Value = 864.44446
$ws.Range("J27").Value = 1374.8334
$ws.Range("K27").Value = 864.44446
$ws.Range("L27").Value = 1374.8334
$ws.Range("M27").Value = -757.44446
$ws.Range("N27").Value = -1588.8334

$ws.Range("H68").Value = 2800
$ws.Range("J68").Value = 3533.3333
$ws.Range("L68").Value = 3533.3333
$ws.Range("N68").Value = -5031.3333

$ws.Range("H71").Value = 2800
$ws.Range("J71").Value = 3533.3333
$ws.Range("L71").Value = 17666.6665
$ws.Range("N71").Value = -25154.6665

$ws.Range("H126").Value = 6069.0645
$ws.Range("I126").Value = 1941.2222
$ws.Range("K126").Value = 5823.6666
$ws.Range("M126").Value = -3353.6666

$ws.Range("H130").Value = 29000
$ws.Range("J130").Value = 29000
$ws.Range("L130").Value = 29000
$ws.Range("N130").Value = -39040

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 13626.571
$ws.Range("I41").Value = 11945
$ws.Range("J41").Value = 14299.2
$ws.Range("K41").Value = 11945
$ws.Range("L41").Value = 14299.2
$ws.Range("M41").Value = -11555
$ws.Range("N41").Value = -15079.2

$ws.Range("H96").Value = 1666
$ws.Range("I96").Value = 1499
$ws.Range("K96").Value = 1499
$ws.Range("M96").Value = -126

$ws.Range("H122").Value = 778.16
$ws.Range("I122").Value = 777.6087
$ws.Range("K122").Value = 2332.8261
$ws.Range("M122").Value = 117.1738999999998

$ws.Range("H136").Value = 3042.7058
$ws.Range("I136").Value = 2918.7144
$ws.Range("K136").Value = 8756.143199999999
$ws.Range("M136").Value = -6206.143199999999
